# Remove the placeholder/guidance paragraph under the "Results Dashboard"
# heading: "(Luca will tackle this. We'll use GitHub pages.)" is no longer
# needed, so delete the whole paragraph (including its paragraph mark).

$d = $word.ActiveDocument

foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*Luca will tackle this*GitHub pages*") {
        $p.Range.Delete()
        break
    }
}
